$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute_options")

$ws.Rows.Item(32).Resize(2).Insert()

$ws.Cells.Item(32, 1).Value = "platform6"
$ws.Cells.Item(32, 2).Value = 6
$ws.Cells.Item(32, 3).Value = "platform"

$ws.Cells.Item(33, 1).Value = "platform7"
$ws.Cells.Item(33, 2).Value = "1-2"
$ws.Cells.Item(33, 3).Value = "platform"

$ws.Range("B32:B33").NumberFormat = "@"
